$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.959.38'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '2.333.71'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.45'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.35'
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.512'
$ws.Range("E7").Value = '  -3.64%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.510'
$ws.Range("E9").Value = '  -3.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.77'
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.26'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0799'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.83'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.86'
$ws.Range("E15").Value = '  +6.05%  '
$ws.Range("D16").Value = '2.352.47'
$ws.Range("E16").Value = '  +2.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.812'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Value = '42.922.42'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.78'
$ws.Range("E19").Value = '  -3.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").Value = '0.0₃0911'
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.82'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.66'
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.38'
$ws.Range("E27").Value = '  +2.48%  '
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("E29").Value = '  +3.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.91'
$ws.Range("E30").Value = '  -4.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.38'
$ws.Range("E31").Value = '  -2.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.22'
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.11'
$ws.Range("E34").Value = '  -2.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.59'
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.60'
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0726'
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("E40").Value = '  -4.56%  '
$ws.Range("E41").Value = '  -3.58%  '
$ws.Range("E42").Value = '  -1.93%  '
$ws.Range("D43").Value = '2.002.94'
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0285'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.78'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.19'
$ws.Range("E46").Value = '  +3.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.93'
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.95'
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").Value = '2.563.12'
$ws.Range("E50").Value = '  +1.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.70'
$ws.Range("E51").Value = '  +2.72%  '
